# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''60.610.33'
$ws.Range('E2').Value = '''  +2.78%  '
$ws.Range('D3').Value = '''2.712.84'
$ws.Range('E3').Value = '''  +2.94%  '
$ws.Range('E4').Value = '''  -0.12%  '
$ws.Range('D5').Value = '''522.36'
$ws.Range('E5').Value = '''  +1.74%  '
$ws.Range('D6').Value = '''146.73'
$ws.Range('E6').Value = '''  +1.99%  '
$ws.Range('D7').Value = '''0.997'
$ws.Range('E7').Value = '''  +0.18%  '
$ws.Range('D8').Value = '''0.576'
$ws.Range('E8').Value = '''  +1.62%  '
$ws.Range('D9').Value = '''2.711.29'
$ws.Range('E9').Value = '''  +1.72%  '
$ws.Range('E10').Value = '''  +2.60%  '
$ws.Range('E11').Value = '''  +0.42%  '
$ws.Range('D12').Value = '''0.342'
$ws.Range('E12').Value = '''  +1.99%  '
$ws.Range('E13').Value = '''  +1.70%  '
$ws.Range('D14').Value = '''3.163.97'
$ws.Range('E14').Value = '''  +2.14%  '
$ws.Range('D15').Value = '''60.782.69'
$ws.Range('E15').Value = '''  +3.14%  '
$ws.Range('D16').Value = '''21.36'
$ws.Range('E16').Value = '''  +1.44%  '
$ws.Range('B17').Value = '''WrappedEther'
$ws.Range('C17').Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '''2.772.57'
$ws.Range('E17').Value = '''  +4.41%  '
$ws.Range('B18').Value = '''ShibaInu'
$ws.Range('C18').Value = '''https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '''0.0000139'
$ws.Range('E18').Value = '''  +1.74%  '
$ws.Range('D19').Value = '''351.78'
$ws.Range('E19').Value = '''  +3.14%  '
$ws.Range('E20').Value = '''  +0.41%  '
$ws.Range('E21').Value = '''  +1.77%  '
$ws.Range('D22').Value = '''6.33'
$ws.Range('E22').Value = '''  +3.91%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '''  +0.07%  '
$ws.Range('D24').Value = '''63.24'
$ws.Range('E24').Value = '''  +3.76%  '
$ws.Range('E25').Value = '''  +0.88%  '
$ws.Range('E26').Value = '''  +5.22%  '
$ws.Range('D27').Value = '''1.01'
$ws.Range('E27').Value = '''  +1.48%  '
$ws.Range('D28').Value = '''0.0₃0817'
$ws.Range('E28').Value = '''  +1.82%  '
$ws.Range('E29').Value = '''  +2.43%  '
$ws.Range('D30').Value = '''6.89'
$ws.Range('E30').Value = '''  +7.88%  '
$ws.Range('E31').Value = '''  +0.15%  '
$ws.Range('E32').Value = '''  +1.75%  '
$ws.Range('D33').Value = '''19.12'
$ws.Range('E33').Value = '''  +1.18%  '
$ws.Range('D34').Value = '''148.92'
$ws.Range('E34').Value = '''  -0.03%  '
$ws.Range('E35').Value = '''  +8.59%  '
$ws.Range('E36').Value = '''  +7.90%  '
$ws.Range('D37').Value = '''0.951'
$ws.Range('E37').Value = '''  -5.97%  '
$ws.Range('D38').Value = '''1.55'
$ws.Range('E38').Value = '''  +10.60%  '
$ws.Range('D39').Value = '''0.884'
$ws.Range('E39').Value = '''  +3.71%  '
$ws.Range('D40').Value = '''36.87'
$ws.Range('E40').Value = '''  +0.84%  '
$ws.Range('E41').Value = '''  +0.91%  '
$ws.Range('D42').Value = '''282.39'
$ws.Range('E42').Value = '''  +0.22%  '
$ws.Range('E43').Value = '''  -0.06%  '
$ws.Range('D44').Value = '''20.04'
$ws.Range('E44').Value = '''  +2.84%  '
$ws.Range('D45').Value = '''0.0989'
$ws.Range('D46').Value = '''0.997'
$ws.Range('E46').Value = '''  +0.25%  '
$ws.Range('D47').Value = '''2.128.21'
$ws.Range('E47').Value = '''  +7.27%  '
$ws.Range('E48').Value = '''  +1.60%  '
$ws.Range('D49').Value = '''4.89'
$ws.Range('E49').Value = '''  +4.03%  '
$ws.Range('D50').Value = '''0.0236'
$ws.Range('E50').Value = '''  +2.88%  '
$ws.Range('B51').Value = '''InjectiveProtocol'
$ws.Range('C51').Value = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '''19.40'
$ws.Range('E51').Value = '''  +5.99%  '
